$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# "Experimental" row (row 7): the Value cell (B7) was empty; fill it in with the
# literal text "true". A direct Range.Value/.Value2 = "true" assignment gets
# auto-coerced to a real Boolean by this host (same as typing it straight into
# a cell), so instead we build the literal text via a formula in a scratch
# cell and paste its *value* back in - that keeps it as plain text, matching
# the source data (Apache POI writes this column as text, never booleans).
$scratch = $ws.Range("D1")
$scratch.Formula = '="true"'
$scratch.Copy()
$target = $ws.Range("B7")
$target.PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()

# "Date" row (row 8): bump the generated timestamp.
$ws.Range("B8").Value = "2025-07-21T12:46:15+00:00"
